$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the stray bold/applied-font style on E3 (back to default/Normal) ---
$ws.Range("E3").Font.Bold = $false

# --- Add three new BOM rows (Main Power connector parts) ---
# Write order matches the original author's entry order so that new shared
# strings land in the same sequence: WM4622-ND, Main Power Male,
# Main Power Female, Main Power Crimp, WM2124-ND, WM20948CT-ND.
$ws.Cells.Item(15, 5).Value = "WM4622-ND"
$ws.Cells.Item(15, 1).Value = "Main Power Male"
$ws.Cells.Item(16, 1).Value = "Main Power Female"
$ws.Cells.Item(17, 1).Value = "Main Power Crimp"
$ws.Cells.Item(16, 5).Value = "WM2124-ND"
$ws.Cells.Item(17, 5).Value = "WM20948CT-ND"

$ws.Cells.Item(15, 2).Value = "Molex"
$ws.Cells.Item(15, 3).Value = 26604040
$ws.Cells.Item(15, 4).Value = "DigiKey"
$ws.Cells.Item(15, 6).Value = 1

$ws.Cells.Item(16, 2).Value = "Molex"
$ws.Cells.Item(16, 3).Value = 9508041
$ws.Cells.Item(16, 4).Value = "DigiKey"
$ws.Cells.Item(16, 6).Value = 1

$ws.Cells.Item(17, 2).Value = "Molex"
$ws.Cells.Item(17, 3).Value = 8500008
$ws.Cells.Item(17, 4).Value = "DigiKey"
$ws.Cells.Item(17, 6).Value = 4

# --- Update the active selection / view to reflect where editing ended ---
$ws.Range("F18").Select()
